# Apply updated Betfair back/lay odds for 2025-12-25 fixtures (rows 2-9)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 3.55
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 2.16
$ws.Range("I2").Value = 2.26
$ws.Range("K2").Value = 3.7
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 1.82
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.31
$ws.Range("T2").Value = 1.8
$ws.Range("U2").Value = 2.04
$ws.Range("V2").Value = 1.79
$ws.Range("W2").Value = 1.33
$ws.Range("X2").Value = 14.5
$ws.Range("Y2").Value = 10.5
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 32
$ws.Range("AB2").Value = 14
$ws.Range("AC2").Value = 7.8
$ws.Range("AD2").Value = 11.5
$ws.Range("AF2").Value = 29
$ws.Range("AG2").Value = 15.5
$ws.Range("AH2").Value = 18.5
$ws.Range("AI2").Value = 44
$ws.Range("AJ2").Value = 85
$ws.Range("AK2").Value = 50
$ws.Range("AL2").Value = 60
$ws.Range("AN2").Value = 60

# Row 3
$ws.Range("F3").Value = 1.97
$ws.Range("G3").Value = 2.06
$ws.Range("H3").Value = 3.95
$ws.Range("I3").Value = 4.3
$ws.Range("K3").Value = 3.95
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.65
$ws.Range("P3").Value = 1.89
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.33
$ws.Range("S3").Value = 3.6
$ws.Range("T3").Value = 1.8
$ws.Range("U3").Value = 2.06
$ws.Range("V3").Value = 1.3
$ws.Range("W3").Value = 1.94
$ws.Range("X3").Value = 14
$ws.Range("Y3").Value = 18
$ws.Range("Z3").Value = 36
$ws.Range("AA3").Value = 95
$ws.Range("AC3").Value = 8.199999999999999
$ws.Range("AD3").Value = 17
$ws.Range("AE3").Value = 95
$ws.Range("AF3").Value = 12
$ws.Range("AG3").Value = 10.5
$ws.Range("AH3").Value = 20
$ws.Range("AK3").Value = 22
$ws.Range("AL3").Value = 40
$ws.Range("AM3").Value = 120
$ws.Range("AN3").Value = 16.5
$ws.Range("AO3").Value = 60

# Row 4
$ws.Range("F4").Value = 2.56
$ws.Range("G4").Value = 2.78
$ws.Range("H4").Value = 2.94
$ws.Range("I4").Value = 3.3
$ws.Range("N4").Value = 3.25
$ws.Range("U4").Value = 1.98
$ws.Range("V4").Value = 1.43
$ws.Range("W4").Value = 1.56
$ws.Range("Y4").Value = 13.5
$ws.Range("Z4").Value = 1000
$ws.Range("AB4").Value = 10
$ws.Range("AD4").Value = 14
$ws.Range("AF4").Value = 20
$ws.Range("AG4").Value = 13
$ws.Range("AI4").Value = 55
$ws.Range("AJ4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("F5").Value = 1.48
$ws.Range("G5").Value = 1.53
$ws.Range("H5").Value = 7.6
$ws.Range("I5").Value = 9.199999999999999
$ws.Range("J5").Value = 4.6
$ws.Range("K5").Value = 5.4
$ws.Range("L5").Value = 1.37
$ws.Range("N5").Value = 4.4
$ws.Range("O5").Value = 1.26
$ws.Range("P5").Value = 2.2
$ws.Range("Q5").Value = 1.74
$ws.Range("R5").Value = 1.46
$ws.Range("S5").Value = 2.96
$ws.Range("T5").Value = 1.96
$ws.Range("U5").Value = 1.92
$ws.Range("W5").Value = 2.88
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 34
$ws.Range("AB5").Value = 9.800000000000001
$ws.Range("AC5").Value = 14
$ws.Range("AD5").Value = 1000
$ws.Range("AF5").Value = 11.5
$ws.Range("AJ5").Value = 15.5
$ws.Range("AK5").Value = 18.5
$ws.Range("AN5").Value = 8.199999999999999

# Row 6
$ws.Range("G6").Value = 3.3
$ws.Range("H6").Value = 2.36
$ws.Range("I6").Value = 2.46
$ws.Range("J6").Value = 3.55
$ws.Range("L6").Value = 1.4
$ws.Range("P6").Value = 1.98
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.37
$ws.Range("S6").Value = 3.55
$ws.Range("T6").Value = 1.74
$ws.Range("V6").Value = 1.68
$ws.Range("X6").Value = 14.5
$ws.Range("Z6").Value = 15.5
$ws.Range("AA6").Value = 34
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 11.5
$ws.Range("AE6").Value = 26
$ws.Range("AF6").Value = 22
$ws.Range("AG6").Value = 14
$ws.Range("AH6").Value = 18
$ws.Range("AI6").Value = 40
$ws.Range("AK6").Value = 38
$ws.Range("AL6").Value = 48
$ws.Range("AM6").Value = 100
$ws.Range("AN6").Value = 34
$ws.Range("AO6").Value = 27

# Row 7
$ws.Range("F7").Value = 1.33
$ws.Range("G7").Value = 1.35
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 12
$ws.Range("J7").Value = 6
$ws.Range("K7").Value = 6.4
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 5.1
$ws.Range("O7").Value = 1.21
$ws.Range("P7").Value = 2.46
$ws.Range("Q7").Value = 1.64
$ws.Range("R7").Value = 1.56
$ws.Range("S7").Value = 2.62
$ws.Range("T7").Value = 1.95
$ws.Range("U7").Value = 1.89
$ws.Range("V7").Value = 1.09
$ws.Range("W7").Value = 3.75
$ws.Range("X7").Value = 24
$ws.Range("Y7").Value = 990
$ws.Range("Z7").Value = 110
$ws.Range("AA7").Value = 470
$ws.Range("AB7").Value = 12
$ws.Range("AC7").Value = 13.5
$ws.Range("AD7").Value = 40
$ws.Range("AE7").Value = 170
$ws.Range("AH7").Value = 29
$ws.Range("AI7").Value = 140
$ws.Range("AJ7").Value = 11
$ws.Range("AK7").Value = 14
$ws.Range("AL7").Value = 36
$ws.Range("AM7").Value = 160
$ws.Range("AN7").Value = 5.2
$ws.Range("AO7").Value = 200

# Row 8
$ws.Range("F8").Value = 1.79
$ws.Range("G8").Value = 1.82
$ws.Range("H8").Value = 6.2
$ws.Range("I8").Value = 7.2
$ws.Range("J8").Value = 3.3
$ws.Range("K8").Value = 3.5
$ws.Range("L8").Value = 1.64
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 2.44
$ws.Range("O8").Value = 1.6
$ws.Range("P8").Value = 1.47
$ws.Range("Q8").Value = 2.84
$ws.Range("R8").Value = 1.16
$ws.Range("S8").Value = 6
$ws.Range("T8").Value = 2.48
$ws.Range("U8").Value = 1.55
$ws.Range("V8").Value = 1.17
$ws.Range("W8").Value = 2.2
$ws.Range("X8").Value = 990
$ws.Range("AB8").Value = 11
$ws.Range("AC8").Value = 990

# Row 9
$ws.Range("G9").Value = 2.26
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 4.2
$ws.Range("N9").Value = 3.75
$ws.Range("P9").Value = 1.96
$ws.Range("Q9").Value = 1.79
$ws.Range("R9").Value = 1.4
$ws.Range("T9").Value = 1.66
$ws.Range("U9").Value = 2.08
$ws.Range("W9").Value = 1.8
